$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    1 = @(74.68892899475621, 65.53000827120312, 74.74188805366863, 44.57478156206889)
    2 = @(82.18830467834425, 86.08869682104951, 89.70325366658092, 57.58147059387897)
    3 = @(77.10845606866206, 78.94736842105263, 85.41844990062607, 55.9817849898298)
    4 = @(61.84084445161264, 69.23076923076923, 71.91099027191366, 47.45809110621789)
    5 = @(89.88957635196672, 59.30065022912257, 80.11025142721257, 48.12989550573491)
    6 = @(90.86068932359662, 50, 66.20471320557093, 30.50546296633821)
    7 = @(69.69683544783204, 73.68421052631578, 81.93583847496714, 41.95455089912883)
    8 = @(57.10337272535252, 46.98837094131527, 56.48564413547995, 30.37860614754947)
    9 = @(68.82335291068271, 60, 66.16596334699783, 44.60839028787305)
}

foreach ($row in $values.Keys) {
    $rowValues = $values[$row]
    $ws.Range("D$row").Value = $rowValues[0]
    $ws.Range("E$row").Value = $rowValues[1]
    $ws.Range("F$row").Value = $rowValues[2]
    $ws.Range("G$row").Value = $rowValues[3]
}
